# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect refreshed scrape output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 11511
$ws1.Range("F9").Value = 4345
$ws1.Range("F14").Value = 2532
$ws1.Range("F16").Value = 130
$ws1.Range("F18").Value = 2214
$ws1.Range("F21").Value = 11298
$ws1.Range("F22").Value = 11207

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 11511
$ws4.Range("F9").Value = 4345
$ws4.Range("F14").Value = 2532
$ws4.Range("F17").Value = 130
$ws4.Range("F19").Value = 2214
$ws4.Range("F22").Value = 11298
$ws4.Range("F23").Value = 11207
